# Apply the edits described by the commit "initial edits to instrument input spreadsheet"
# to the instrument_definition_inputs workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update header row (C1/D1 labels) -------------------------------------
# "Nwafer" -> "Nchan", "det_spacing_freq" -> "Nhorns"
$ws.Range("C1").Value = "Nchan"
$ws.Range("D1").Value = "Nhorns"

# --- Update the data table (rows 2-9) --------------------------------------
$ws.Range("A2").Value = 92.453749999999999
$ws.Range("B2").Value = 92.546250000000001
$ws.Range("C2").Value = 378
$ws.Range("D2").Value = 500

$ws.Range("A3").Value = 149.92500000000001
$ws.Range("B3").Value = 150.07499999999999
$ws.Range("C3").Value = 266
$ws.Range("D3").Value = 500

$ws.Range("A4").Value = 219.89
$ws.Range("B4").Value = 220.11
$ws.Range("C4").Value = 227
$ws.Range("D4").Value = 500

$ws.Range("A5").Value = 277.36124999999998
$ws.Range("B5").Value = 277.63875000000002
$ws.Range("C5").Value = 234
$ws.Range("D5").Value = 500

$ws.Range("A6").Value = 349.82499999999999
$ws.Range("B6").Value = 350.17500000000001
$ws.Range("C6").Value = 85
$ws.Range("D6").Value = 500

$ws.Range("A7").Value = 404.79750000000001
$ws.Range("B7").Value = 405.20249999999999
$ws.Range("C7").Value = 74
$ws.Range("D7").Value = 500

$ws.Range("A8").Value = 667.16624999999999
$ws.Range("B8").Value = 667.83375000000001
$ws.Range("C8").Value = 127
$ws.Range("D8").Value = 500

$ws.Range("A9").Value = 872.06375000000003
$ws.Range("B9").Value = 872.93624999999997
$ws.Range("C9").Value = 63
$ws.Range("D9").Value = 500

# --- Remove the old summary row (row 10 with the NaN labels & SUM formula) -
$ws.Rows.Item(10).Delete()

# --- Update the selected cell to match the author's final selection --------
$ws.Range("E9").Select()
